$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prices")

$row = 97

$values = @("2025-06-06", "35.5", "35.21", "0.94", "0.248", "0.09", "5,512", "8,253", "8,303", "7.1975")

for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $values[$col - 1]
    $cell.Style = "Normal"
}
